$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the last data row (134) down to row 135 by copying it (this preserves
# styles/number formats exactly, unlike Rows.Insert which recomputes styles).
$ws.Range("A134:AC134").Copy($ws.Range("A135"))

# --- New row 134: a new match result ---
$ws.Range("A134").Value2 = 132
$ws.Range("B134").Value2 = 7011626
$ws.Range("C134").Value = "Azerbaijan Premier League"
$ws.Range("D134").Value = "Azerbaijan Premier League"
$ws.Range("E134").Value2 = 45368.375
$ws.Range("F134").Value = "PFK Turan Tovuz"
$ws.Range("G134").Value = "FK Qarabag"
$ws.Range("H134").Value2 = 1
$ws.Range("I134").Value2 = 3
$ws.Range("J134").Value = "A"
$ws.Range("K134").Value2 = 4.333
$ws.Range("L134").Value2 = 4
$ws.Range("M134").Value2 = 1.571
$ws.Range("N134").Value2 = 4.5
$ws.Range("O134").Value2 = 4
$ws.Range("P134").Value2 = 1.55
$ws.Range("Q134").Value2 = 1
$ws.Range("R134").Value2 = 1.8
$ws.Range("S134").Value2 = 2
$ws.Range("T134").Value2 = 2.75
$ws.Range("U134").Value2 = 1.875
$ws.Range("V134").Value2 = 1.925
$ws.Range("W134").Value2 = -1
$ws.Range("X134").Value2 = -1
$ws.Range("Y134").Value2 = 0.55
$ws.Range("Z134").Value2 = -1
$ws.Range("AA134").Value2 = 1
$ws.Range("AB134").Value2 = 0.875
$ws.Range("AC134").Value2 = -1

# --- Row 135 (previously row 134): update with final match odds/result ---
$ws.Range("A135").Value2 = 133
$ws.Range("H135").Value2 = 1
$ws.Range("I135").Value2 = 2
$ws.Range("J135").Value = "A"
$ws.Range("K135").Value2 = 2.9
$ws.Range("L135").Value2 = 3.75
$ws.Range("M135").Value2 = 2
$ws.Range("N135").Value2 = 2.75
$ws.Range("O135").Value2 = 3.5
$ws.Range("P135").Value2 = 2.15
$ws.Range("Q135").Value2 = 0.25
$ws.Range("R135").Value2 = 1.85
$ws.Range("S135").Value2 = 1.95
$ws.Range("T135").Value2 = 2.25
$ws.Range("U135").Value2 = 1.95
$ws.Range("V135").Value2 = 1.85
$ws.Range("W135").Value2 = -1
$ws.Range("X135").Value2 = -1
$ws.Range("Y135").Value2 = 1.15
$ws.Range("Z135").Value2 = -1
$ws.Range("AA135").Value2 = 0.95
$ws.Range("AB135").Value2 = 0.95
$ws.Range("AC135").Value2 = -1
